# Update cryptocurrency price (D) and 1h volume/change (E) columns
# to reflect the latest scrape, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [pscustomobject]@{ Row = 2; Price = '22.380.08'; Volume = '  -0.14%  ' }
    [pscustomobject]@{ Row = 3; Price = '1.568.05'; Volume = '  -0.27%  ' }
    [pscustomobject]@{ Row = 4; Price = $null; Volume = '  +0.14%  ' }
    [pscustomobject]@{ Row = 5; Price = $null; Volume = '  +0.11%  ' }
    [pscustomobject]@{ Row = 6; Price = '291.06'; Volume = '  +0.45%  ' }
    [pscustomobject]@{ Row = 7; Price = '0.3781'; Volume = '  +2.95%  ' }
    [pscustomobject]@{ Row = 8; Price = '49.10'; Volume = '  -0.42%  ' }
    [pscustomobject]@{ Row = 9; Price = '0.3403'; Volume = '  -0.01%  ' }
    [pscustomobject]@{ Row = 10; Price = '0.07613'; Volume = '  -0.46%  ' }
    [pscustomobject]@{ Row = 11; Price = '1.141'; Volume = '  -2.65%  ' }
    [pscustomobject]@{ Row = 12; Price = '1.002'; Volume = '  +0.11%  ' }
    [pscustomobject]@{ Row = 13; Price = '21.10'; Volume = '  -0.83%  ' }
    [pscustomobject]@{ Row = 14; Price = '5.989'; Volume = '  -1.40%  ' }
    [pscustomobject]@{ Row = 15; Price = '6.932'; Volume = $null }
    [pscustomobject]@{ Row = 16; Price = '1.567.54'; Volume = '  +0.01%  ' }
    [pscustomobject]@{ Row = 17; Price = '0.00001133'; Volume = '  -0.04%  ' }
    [pscustomobject]@{ Row = 18; Price = '89.97'; Volume = '  +0.08%  ' }
    [pscustomobject]@{ Row = 19; Price = '0.06737'; Volume = $null }
    [pscustomobject]@{ Row = 20; Price = $null; Volume = '  +0.15%  ' }
    [pscustomobject]@{ Row = 21; Price = '16.64'; Volume = '  +0.32%  ' }
    [pscustomobject]@{ Row = 22; Price = '6.205'; Volume = '  -0.88%  ' }
    [pscustomobject]@{ Row = 23; Price = '11.96'; Volume = '  -0.48%  ' }
    [pscustomobject]@{ Row = 24; Price = '22.381.79'; Volume = '  -0.10%  ' }
    [pscustomobject]@{ Row = 25; Price = '2.407'; Volume = '  +1.90%  ' }
    [pscustomobject]@{ Row = 26; Price = '2.703'; Volume = '  -7.14%  ' }
    [pscustomobject]@{ Row = 27; Price = '20.17'; Volume = '  +0.59%  ' }
    [pscustomobject]@{ Row = 28; Price = '147.14'; Volume = '  +0.57%  ' }
    [pscustomobject]@{ Row = 29; Price = '5.019'; Volume = '  +0.59%  ' }
    [pscustomobject]@{ Row = 30; Price = '126.20'; Volume = '  +0.35%  ' }
    [pscustomobject]@{ Row = 31; Price = '1.739.07'; Volume = '  -0.16%  ' }
    [pscustomobject]@{ Row = 32; Price = '2.018'; Volume = '  +0.00%  ' }
    [pscustomobject]@{ Row = 33; Price = '6.112'; Volume = '  -2.28%  ' }
    [pscustomobject]@{ Row = 34; Price = '0.9965'; Volume = '  -2.52%  ' }
    [pscustomobject]@{ Row = 35; Price = '10.14'; Volume = '  +0.08%  ' }
    [pscustomobject]@{ Row = 36; Price = $null; Volume = '  +9.17%  ' }
    [pscustomobject]@{ Row = 37; Price = $null; Volume = '  +0.58%  ' }
    [pscustomobject]@{ Row = 38; Price = '0.02518'; Volume = '  -1.12%  ' }
    [pscustomobject]@{ Row = 39; Price = '0.2299'; Volume = '  -1.13%  ' }
    [pscustomobject]@{ Row = 40; Price = '0.06496'; Volume = '  +0.02%  ' }
    [pscustomobject]@{ Row = 41; Price = '5.414'; Volume = '  -2.08%  ' }
    [pscustomobject]@{ Row = 42; Price = '11.38'; Volume = '  -3.20%  ' }
    [pscustomobject]@{ Row = 43; Price = '0.6337'; Volume = '  -0.46%  ' }
    [pscustomobject]@{ Row = 44; Price = $null; Volume = '  +0.15%  ' }
    [pscustomobject]@{ Row = 45; Price = '13.98'; Volume = '  -1.59%  ' }
    [pscustomobject]@{ Row = 46; Price = '3.805'; Volume = '  +1.26%  ' }
    [pscustomobject]@{ Row = 47; Price = '0.5941'; Volume = '  -1.06%  ' }
    [pscustomobject]@{ Row = 48; Price = '2.087'; Volume = '  -1.20%  ' }
    [pscustomobject]@{ Row = 49; Price = '1.252'; Volume = '  -0.48%  ' }
    [pscustomobject]@{ Row = 50; Price = '124.52'; Volume = '  -0.18%  ' }
    [pscustomobject]@{ Row = 51; Price = '0.07325'; Volume = '  +0.45%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)   # column D = Price
        $priceCell.NumberFormat = "@"            # keep as text, matching source data
        $priceCell.Value = $u.Price
    }
    if ($null -ne $u.Volume) {
        $ws.Cells.Item($u.Row, 5).Value = $u.Volume   # column E = Volume(1h)
    }
}

Write-Host "Updated $($updates.Count) rows in cryptos sheet"
